$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.838.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '''1.856.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("D4").Value = '''0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''304.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").Value = '''0.9999'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  -1.85%  '
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").Value = '''0.07154'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").Value = '''0.8909'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = '''20.65'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '''1.859.84'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").Value = '''0.07440'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.90%  '
$ws.Range("D14").Value = '''92.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("D15").Value = '''5.219'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '''1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '''0.000008509'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '''26.878.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.23%  '
$ws.Range("D21").Value = '''5.015'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").Value = '''2.095.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").Value = '''6.433'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").Value = '''1.793'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.89%  '
$ws.Range("D27").Value = '''17.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").Value = '''2.057'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("D29").Value = '''112.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("D31").Value = '''4.653'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").Value = '''0.09214'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.35%  '
$ws.Range("D33").Value = '''0.05077'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("D34").Value = '''2.985'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("D35").Value = '''0.7422'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.69%  '
$ws.Range("D36").Value = '''1.144'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("D37").Value = '''3.234'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.81%  '
$ws.Range("D38").Value = '''2.519'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("D39").Value = '''0.01986'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.85%  '
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").Value = '''0.5321'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").Value = '''119.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.45%  '
$ws.Range("D43").Value = '''6.467'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("D44").Value = '''8.358'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").Value = '''0.1454'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").Value = '''1.0000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").Value = '''0.4631'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").Value = '''10.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").Value = '''1.557'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").Value = '''36.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '''62.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.60%  '
